$d = $word.ActiveDocument

# --- 1. Lambda paragraph: drop the _GoBack bookmark, append a new sentence ---
$bm = $d.Bookmarks("_GoBack")
$bmStart = $bm.Start
$bm.Delete()
$ip = $d.Range($bmStart, $bmStart)
$ip.InsertAfter(" Tuple type as an argument to a lambda expression.")

# --- 2. Fill in the (previously empty) last paragraph with the "const" note ---
$pConst = $d.Paragraphs.Last
$rConst = $pConst.Range
$rConst.Collapse(0)
$rConst.InsertAfter("The static is not allowed in a constant declaration. Readonly vs constant, const can only be initialized at the declaration.")

# --- 3. New paragraph: readonly ---
$rConst.Collapse(0)
$rConst.InsertParagraphAfter()
$pReadonly = $d.Paragraphs.Last
$rReadonly = $pReadonly.Range
$rReadonly.Collapse(0)
$rReadonly.InsertAfter("Readonly can be initialized in the declaration, in the constructor. Readonly struct(struct is immutable), ref readonly(returned reference cannot be modified)")

# --- 4. New paragraph: events ---
$rReadonly.Collapse(0)
$rReadonly.InsertParagraphAfter()
$pEvents = $d.Paragraphs.Last
$rEvents = $pEvents.Range
$rEvents.Collapse(0)
$rEvents.InsertAfter("Events: publisher(when) and receiver(what action), e.g. button click")

# --- 5. New paragraph: extern (+ bookmark relocated here) ---
$rEvents.Collapse(0)
$rEvents.InsertParagraphAfter()
$pExtern = $d.Paragraphs.Last
$rExtern = $pExtern.Range
$rExtern.Collapse(0)
$rExtern.InsertAfter("Extern: combine with DllImport, must be declared as static")

$rExtern.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rExtern)
